$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CVS PHARMACY"
$ws.Range("B2").Value = "CVS/PHARMACY"
$ws.Range("A3").Value = "WAL-MART"
$ws.Range("B3").Value = "WALMART"
$ws.Range("A1").Value = "actualName"
$ws.Range("B1").Value = "expectedName"

$ws.Range("A6").Select()
